# Add the "Статистика по городам" worksheet after the existing sheet,
# populate it with city salary / vacancy-share data and formatting that
# matches the target workbook.

$wb = $excel.ActiveWorkbook

# --- create the new sheet, positioned after all existing sheets ---
$lastIndex = $wb.Worksheets.Count
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$ws.Name = "Статистика по городам"
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# --- column widths (stored OOXML width = ColumnWidth + 5/6) ---
$ws.Columns.Item(1).ColumnWidth = 16 - 5/6   # A
$ws.Columns.Item(2).ColumnWidth = 16 - 5/6   # B
$ws.Columns.Item(4).ColumnWidth = 16 - 5/6   # D
$ws.Columns.Item(5).ColumnWidth = 21 - 5/6   # E

# --- header row ---
$ws.Range("A1").Value = "Город"
$ws.Range("B1").Value = "Уровень зарплат"
$ws.Range("D1").Value = "Город"
$ws.Range("E1").Value = "Доля вакансий"

# Headers use the bold/bordered style already present on sheet1's header row.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("D1:E1").PasteSpecial(-4122)

# --- data rows ---
$rows = @(
    @("Москва", 40201, "Москва", 0.334026018899322),
    @("Санкт-Петербург", 38461, "Санкт-Петербург", 0.1014048968812217),
    @("Новосибирск", 34765, "Минск", 0.04061388562982566),
    @("Екатеринбург", 32600, "Киев", 0.03216914223228148),
    @("Краснодар", 31798, "Новосибирск", 0.02363188261188267),
    @("Казань", 29840, "Нижний Новгород", 0.02144437701366958),
    @("Ростов-на-Дону", 28090, "Екатеринбург", 0.01972705996450027),
    @("Пермь", 27517, "Алматы", 0.01892116268106615),
    @("Самара", 27097, "Воронеж", 0.01859600619490592),
    @("Нижний Новгород", 26376, "Казань", 0.01856606725806701)
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $r = $r + 1
}

# Data cells A2:B11 and D2:D11 reuse the plain bordered style from sheet1.
$ws1.Range("A2").Copy()
$ws.Range("A2:B11").PasteSpecial(-4122)
$ws.Range("D2:D11").PasteSpecial(-4122)

# Column E holds a percentage share -> bordered style + percentage number format.
$ws.Range("E2:E11").PasteSpecial(-4122)
$ws.Range("E2:E11").NumberFormat = "0.00%"

$ws.Range("A1").Select()

# Restore the original sheet as the active tab (the diff keeps activeTab="0").
$ws1.Activate()
